# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" / "Valor Mora" block (rows 16-25, columns E/F) is
# reordered: the periods that used to run ascending (2008 .. 2105) are
# flipped to run descending (2105 .. 2008), and the "Valor Mora" figure
# that belonged to period 2105 (28090, as opposed to the common 35112)
# now travels with it to the top row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Periodo Mora") - reverse the ten-row block.
$ws.Range("E16").Value = "2105"
$ws.Range("E17").Value = "2104"
$ws.Range("E18").Value = "2103"
$ws.Range("E19").Value = "2102"
$ws.Range("E20").Value = "2101"
$ws.Range("E21").Value = "2012"
$ws.Range("E22").Value = "2011"
$ws.Range("E23").Value = "2010"
$ws.Range("E24").Value = "2009"
$ws.Range("E25").Value = "2008"

# Column F ("Valor Mora") - the non-default amount follows period 2105
# from the last row of the block up to the first row.
$ws.Range("F16").Value = 28090
$ws.Range("F25").Value = 35112
